$d = $word.ActiveDocument

# Locate the paragraph that contains "${Footer.FooterChild.nestedData}"
# (built from runs: "${" + "Footer.FooterChild.nestedData" + "}")
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "`${Footer.FooterChild.nestedData}") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph '`${Footer.FooterChild.nestedData}'"
}

$targetParagraph = $d.Paragraphs.Item($targetIndex)

# Insert a new empty paragraph right after it
$targetParagraph.Range.InsertParagraphAfter()

# The newly created paragraph is the one right after the target paragraph
$newParagraph = $d.Paragraphs.Item($targetIndex + 1)

# Build the new paragraph content, mirroring the existing
# "${Footer.FooterChild.nestedData}" paragraph's run/proofErr layout, but
# splitting "Footer.FooterChild." and "nestedData.withPoint" into two runs
# and appending ".withPoint" to the parameter name.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>${</w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>Footer.FooterChild.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>nestedData.withPoint</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>}</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newParagraph.Range.InsertXML($xml)
